$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34, column A: the phone number was stored as text ("71277620");
# re-entering it turns it into a genuine number, matching upstream.
$ws.Range("A34").Value = 71277620

# New row 35 - payment 71277620 (Cash) 2025-08-18T17:10:20
# Keep the phone number (column A) as text, same as the rest of the sheet.
$a35 = $ws.Range("A35")
$a35.NumberFormat = "@"
$a35.Value = "71277620"
$a35.Style = "Normal"

$ws.Range("C35").Value = "Cash"
$ws.Range("D35").Value = "2025-08-18T17:10:20"
$ws.Range("E35").Value = 76
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 76
